$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ASSESSMENTS")
$ws2 = $wb.Worksheets.Item("INDUSTRIES")

# Rename "industry" labels to "site" labels on the INDUSTRIES sheet.
$ws2.Range("B3").Value = "Site "
$ws2.Range("A1").Value = "SITE"
$ws2.Range("E1").Value = "SUB-SUPPLIERS"

# Update selection on ASSESSMENTS sheet (no longer the active tab).
[void]$ws1.Range("B12").Select()

# Make INDUSTRIES the active sheet/tab, with its own new selection.
[void]$ws2.Activate()
[void]$ws2.Range("E1:J1").Select()
